$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last existing row (63) had a special "date only" number format since
# it was the final/most-recent row. A new row (64) is being appended with
# today's data, so row 63's date cell reverts to the standard datetime
# format used throughout the rest of the column, and row 64's date cell
# takes on the special "last row" format.
$ws.Range("A63").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A64").Value = 45651
$ws.Range("B64").Value = 151
$ws.Range("C64").Value = 140
$ws.Range("D64").Value = 148

$ws.Range("A64").NumberFormat = "YYYY-MM-DD"
